$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Summary")

# Row 2
$ws1.Range("B2").Value = 45995.51041666666
$ws1.Range("C2").Value = "Амур"
$ws1.Range("D2").Value = "Локомотив"
$ws1.Range("E2").Value = "Амур – Локомотив"
$ws1.Range("F2").Value = 897848
$ws1.Range("G2").Value = "https://text.khl.ru/text/897848.html"
$ws1.Range("H2").Value = 1.09375
$ws1.Range("I2").Value = 1.323529
$ws1.Range("J2").Value = 2.417279
$ws1.Range("K2").Value = 21.410046
$ws1.Range("L2").Value = 27.038132
$ws1.Range("M2").Value = 48.448178
$ws1.Range("N2").Value = 0.30718
$ws1.Range("O2").Value = 0.261547
$ws1.Range("P2").Value = 0.431273
$ws1.Range("Q2").Value = 3.255420274757471
$ws1.Range("R2").Value = 3.823404588850188
$ws1.Range("S2").Value = 2.318716914808022
$ws1.Range("T2").Value = 30.718
$ws1.Range("U2").Value = 26.1547
$ws1.Range("V2").Value = 43.12730000000001
$ws1.Range("W2").Value = 0.878381
$ws1.Range("X2").Value = 0.121619
$ws1.Range("Y2").Value = 8.222399460610594
$ws1.Range("Z2").Value = 0.95146
$ws1.Range("AA2").Value = 0.04854
$ws1.Range("AB2").Value = 20.60156571899464
$ws1.Range("AC2").Value = 0.983045
$ws1.Range("AD2").Value = 0.016955
$ws1.Range("AE2").Value = 58.97965202005307
$ws1.Range("AF2").Value = 0.324304
$ws1.Range("AG2").Value = 0.675696
$ws1.Range("AH2").Value = 3.083526567664907
$ws1.Range("AI2").Value = 0.112789
$ws1.Range("AJ2").Value = 0.887211
$ws1.Range("AK2").Value = 8.866112830151877
$ws1.Range("AL2").Value = 0.418231
$ws1.Range("AM2").Value = 0.581769
$ws1.Range("AN2").Value = 2.391023142712998
$ws1.Range("AO2").Value = 0.173621
$ws1.Range("AP2").Value = 0.826379
$ws1.Range("AQ2").Value = 5.759671929086919
$ws1.Range("AR2").Value = 0.792876
$ws1.Range("AS2").Value = 1.261231264409567
$ws1.Range("AT2").Value = 0.875358
$ws1.Range("AU2").Value = 1.142389742254026

# Row 3
$ws1.Range("B3").Value = 45995.52083333334
$ws1.Range("C3").Value = "Адмирал"
$ws1.Range("D3").Value = "СКА"
$ws1.Range("E3").Value = "Адмирал – СКА"
$ws1.Range("F3").Value = 897849
$ws1.Range("G3").Value = "https://text.khl.ru/text/897849.html"
$ws1.Range("H3").Value = 2.285128
$ws1.Range("I3").Value = 4
$ws1.Range("J3").Value = 6.285128
$ws1.Range("K3").Value = 30.885829
$ws1.Range("L3").Value = 34.182529
$ws1.Range("M3").Value = 65.068358
$ws1.Range("N3").Value = 0.37997
$ws1.Range("O3").Value = 0.157033
$ws1.Range("P3").Value = 0.461572
$ws1.Range("Q3").Value = 2.631786720004211
$ws1.Range("R3").Value = 6.368088236230601
$ws1.Range("S3").Value = 2.166509233662354
$ws1.Range("T3").Value = 37.997
$ws1.Range("U3").Value = 15.7033
$ws1.Range("V3").Value = 46.1572
$ws1.Range("W3").Value = 0.206735
$ws1.Range("X3").Value = 0.791839
$ws1.Range("Y3").Value = 1.262882985051254
$ws1.Range("Z3").Value = 0.346762
$ws1.Range("AA3").Value = 0.651813
$ws1.Range("AB3").Value = 1.534182349845738
$ws1.Range("AC3").Value = 0.5021
$ws1.Range("AD3").Value = 0.496474
$ws1.Range("AE3").Value = 2.014204167791264
$ws1.Range("AF3").Value = 0.827826
$ws1.Range("AG3").Value = 0.172174
$ws1.Range("AH3").Value = 1.207983320166315
$ws1.Range("AI3").Value = 0.618542
$ws1.Range("AJ3").Value = 0.381458
$ws1.Range("AK3").Value = 1.616705090357647
$ws1.Range("AL3").Value = 0.860208
$ws1.Range("AM3").Value = 0.139792
$ws1.Range("AN3").Value = 1.162509532578167
$ws1.Range("AO3").Value = 0.672358
$ws1.Range("AP3").Value = 0.327642
$ws1.Range("AQ3").Value = 1.487302895183816
$ws1.Range("AR3").Value = 0.687709
$ws1.Range("AS3").Value = 1.454103407109693
$ws1.Range("AT3").Value = 0.75752
$ws1.Range("AU3").Value = 1.320097159150914

# Row 4
$ws1.Range("B4").Value = 45995.6875
$ws1.Range("C4").Value = "Авангард"
$ws1.Range("D4").Value = "ХК Сочи"
$ws1.Range("E4").Value = "Авангард – ХК Сочи"
$ws1.Range("F4").Value = 897847
$ws1.Range("G4").Value = "https://text.khl.ru/text/897847.html"
$ws1.Range("H4").Value = 4.5
$ws1.Range("I4").Value = 1
$ws1.Range("J4").Value = 5.5
$ws1.Range("K4").Value = 40.265219
$ws1.Range("L4").Value = 23.307977
$ws1.Range("M4").Value = 63.573196
$ws1.Range("N4").Value = 0.896996
$ws1.Range("O4").Value = 0.055477
$ws1.Range("P4").Value = 0.036569
$ws1.Range("Q4").Value = 1.114832173164652
$ws1.Range("R4").Value = 18.02548804008869
$ws1.Range("S4").Value = 27.34556591648665
$ws1.Range("T4").Value = 89.6996
$ws1.Range("U4").Value = 5.5477
$ws1.Range("V4").Value = 3.6569
$ws1.Range("W4").Value = 0.279957
$ws1.Range("X4").Value = 0.709085
$ws1.Range("Y4").Value = 1.410268162491098
$ws1.Range("Z4").Value = 0.43954
$ws1.Range("AA4").Value = 0.549501
$ws1.Range("AB4").Value = 1.819832902942852
$ws1.Range("AC4").Value = 0.600144
$ws1.Range("AD4").Value = 0.388898
$ws1.Range("AE4").Value = 2.571368327942031
$ws1.Range("AF4").Value = 0.953672
$ws1.Range("AG4").Value = 0.046328
$ws1.Range("AH4").Value = 1.048578546921793
$ws1.Range("AI4").Value = 0.860844
$ws1.Range("AJ4").Value = 0.139156
$ws1.Range("AK4").Value = 1.161650659120584
$ws1.Range("AL4").Value = 0.338198
$ws1.Range("AM4").Value = 0.661802
$ws1.Range("AN4").Value = 2.956847763736036
$ws1.Range("AO4").Value = 0.121009
$ws1.Range("AP4").Value = 0.878991
$ws1.Range("AQ4").Value = 8.263848143526515
$ws1.Range("AR4").Value = 0.977066
$ws1.Range("AS4").Value = 1.023472314050433
$ws1.Range("AT4").Value = 0.19097
$ws1.Range("AU4").Value = 5.236424569304079

$ws2 = $wb.Worksheets.Item("Cards_telegram")

# Row 2
$ws2.Range("A2").Value = 45995.51041666666
$ws2.Range("B2").Value = "Амур – Локомотив"
$ws2.Range("C2").Value = "КХЛ • Регулярный чемпионат • 04.12.2025`n`nАмур – Локомотив`n`nОжидания модели (60’):`n• Голы: λ_total ≈ 2.59 (1.16 : 1.43)`n• Броски: SOG λ ≈ 48 (21 : 27)`n`nИсход (60’), честные кф:`n• П1: 30.7%  (Kмод 3.26)`n• Х:  26.2%  (Kмод 3.82)`n• П2: 43.1%  (Kмод 2.32)`n`nТоталы голов:`n• ТМ 4.5: 87.8%  (Kмод 1.14)`n• ТБ 4.5: 12.2%  (Kмод 8.22)`n`n• ТМ 5.5: 95.1%  (Kмод 1.05)`n• ТБ 5.5: 4.9%  (Kмод 20.60)`n`n• ТМ 6.5: 98.3%  (Kмод 1.02)`n• ТБ 6.5: 1.7%  (Kмод 58.98)`n`nИндивидуальные тоталы:`n• Амур ИТБ 1.5: 32.4% (Kмод 3.08)`n• Амур ИТБ 2.5: 11.3% (Kмод 8.87)`n• Локомотив ИТБ 1.5: 41.8% (Kмод 2.39)`n• Локомотив ИТБ 2.5: 17.4% (Kмод 5.76)`n`nФора +1.5:`n• Амур +1.5: 79.3% (Kмод 1.26)`n• Локомотив +1.5: 87.5% (Kмод 1.14)"

# Row 3
$ws2.Range("A3").Value = 45995.52083333334
$ws2.Range("B3").Value = "Адмирал – СКА"
$ws2.Range("C3").Value = "КХЛ • Регулярный чемпионат • 04.12.2025`n`nАдмирал – СКА`n`nОжидания модели (60’):`n• Голы: λ_total ≈ 6.66 (3.19 : 3.46)`n• Броски: SOG λ ≈ 65 (31 : 34)`n`nИсход (60’), честные кф:`n• П1: 38.0%  (Kмод 2.63)`n• Х:  15.7%  (Kмод 6.37)`n• П2: 46.2%  (Kмод 2.17)`n`nТоталы голов:`n• ТМ 4.5: 20.7%  (Kмод 4.84)`n• ТБ 4.5: 79.2%  (Kмод 1.26)`n`n• ТМ 5.5: 34.7%  (Kмод 2.88)`n• ТБ 5.5: 65.2%  (Kмод 1.53)`n`n• ТМ 6.5: 50.2%  (Kмод 1.99)`n• ТБ 6.5: 49.6%  (Kмод 2.01)`n`nИндивидуальные тоталы:`n• Адмирал ИТБ 1.5: 82.8% (Kмод 1.21)`n• Адмирал ИТБ 2.5: 61.9% (Kмод 1.62)`n• СКА ИТБ 1.5: 86.0% (Kмод 1.16)`n• СКА ИТБ 2.5: 67.2% (Kмод 1.49)`n`nФора +1.5:`n• Адмирал +1.5: 68.8% (Kмод 1.45)`n• СКА +1.5: 75.8% (Kмод 1.32)"

# Row 4
$ws2.Range("A4").Value = 45995.6875
$ws2.Range("B4").Value = "Авангард – ХК Сочи"
$ws2.Range("C4").Value = "КХЛ • Регулярный чемпионат • 04.12.2025`n`nАвангард – ХК Сочи`n`nОжидания модели (60’):`n• Голы: λ_total ≈ 6.04 (4.84 : 1.20)`n• Броски: SOG λ ≈ 64 (40 : 23)`n`nИсход (60’), честные кф:`n• П1: 89.7%  (Kмод 1.11)`n• Х:  5.5%  (Kмод 18.03)`n• П2: 3.7%  (Kмод 27.35)`n`nТоталы голов:`n• ТМ 4.5: 28.0%  (Kмод 3.57)`n• ТБ 4.5: 70.9%  (Kмод 1.41)`n`n• ТМ 5.5: 44.0%  (Kмод 2.28)`n• ТБ 5.5: 55.0%  (Kмод 1.82)`n`n• ТМ 6.5: 60.0%  (Kмод 1.67)`n• ТБ 6.5: 38.9%  (Kмод 2.57)`n`nИндивидуальные тоталы:`n• Авангард ИТБ 1.5: 95.4% (Kмод 1.05)`n• Авангард ИТБ 2.5: 86.1% (Kмод 1.16)`n• ХК Сочи ИТБ 1.5: 33.8% (Kмод 2.96)`n• ХК Сочи ИТБ 2.5: 12.1% (Kмод 8.26)`n`nФора +1.5:`n• Авангард +1.5: 97.7% (Kмод 1.02)`n• ХК Сочи +1.5: 19.1% (Kмод 5.24)"

